# Applies:
#  1) Rename sheet "dados_de_geradores" -> "dados_dos_geradores"
#     Rename sheet "Página4"            -> "dados_dos_transformadores"
#  2) On "dados_de_linha": extend the little styled block in columns A:B
#     down to rows 30-31 (same style as the existing A32:B38 block), add a
#     new centred/wrapped block in columns D:E for rows 34-39, and remove
#     the old A:B styling on rows 37-38 (replaced by the D:E block).

$wb = $excel.ActiveWorkbook

# --- 1) sheet renames ---------------------------------------------------
$wsGeradores = $wb.Worksheets.Item("dados_de_geradores")
$wsGeradores.Name = "dados_dos_geradores"

$wsPagina4 = $wb.Worksheets.Item("Página4")
$wsPagina4.Name = "dados_dos_transformadores"

# --- 2) "dados_de_linha" row/style tweaks -------------------------------
$ws = $wb.Worksheets.Item("dados_de_linha")

# Extend the existing A/B formatting block upward into rows 30-31, reusing
# the same style already applied to A32:B38 (copy formats only).
$ws.Range("A32:B32").Copy()
$ws.Range("A30:B31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Create the new centred/word-wrapped style (no borders, Times New Roman)
# on D34, matching the rest of the sheet's "data" font.
$seed = $ws.Range("D34")
$seed.Font.Name = "Times New Roman"
$seed.Font.Size = 10
$seed.Font.Color = 0
$seed.HorizontalAlignment = -4108
$seed.WrapText = $true
$seed.ShrinkToFit = $false

# Propagate that new style across D34:E39.
$seed.Copy()
$ws.Range("D34:E39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 37-38 lose their old A/B styling entirely (back to "no cell").
$ws.Range("A37:B38").Clear()
